$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stats (runs, balls, fours, sixes) for row 2 and row 4 were swapped
# between two different matches. Swap columns C:F between row 2 and row 4.
$cols = @("C", "D", "E", "F")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range("$col" + "2")
    $cellRow4 = $ws.Range("$col" + "4")

    $row2Value = $cellRow2.Value2
    $row4Value = $cellRow4.Value2

    $cellRow2.Value2 = $row4Value
    $cellRow4.Value2 = $row2Value
}
